$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: add Responsable "Agustina" with 100% progress
$ws.Range("B9").Value = "Agustina"
$ws.Range("C9").Value = 1
$ws.Range("C9").NumberFormat = "0%"

# Row 28: add Responsable "Lucas" with 100% progress
$ws.Range("B28").Value = "Lucas"
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = "0%"

# Row 30: "en proceso" -> 100% (done)
$ws.Range("C30").Value = 1
$ws.Range("C30").NumberFormat = "0%"

# Row 35: "en proceso" -> 100% (done)
$ws.Range("C35").Value = 1
$ws.Range("C35").NumberFormat = "0%"

# Row 36: "en proceso" -> 100% (done)
$ws.Range("C36").Value = 1
$ws.Range("C36").NumberFormat = "0%"

# New row 41: new task reported, assigned to Agustina, 100% done
$ws.Range("A41").Value = "en los reportes faltan los decimales"
$ws.Range("B41").Value = "Agustina"
$ws.Range("C41").Value = 1
$ws.Range("C41").NumberFormat = "0%"

# Update the visible selection to C10 (also drops the stale topLeftCell scroll anchor)
$ws.Range("C10").Select()
